$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.300.80'
$ws.Range("E2").Value = '  +1.10%  '

$ws.Range("D3").Value = '2.478.16'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '563.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '163.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("D9").Value = '2.476.99'
$ws.Range("E9").Value = '  +0.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.153'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.166'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.85%  '

$ws.Range("B12").Value = 'Toncoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.91'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.88%  '

$ws.Range("B13").Value = 'Cardano'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.334'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.60%  '

$ws.Range("D14").Value = '69.159.38'
$ws.Range("E14").Value = '  +1.11%  '

$ws.Range("E15").Value = '  +0.02%  '

$ws.Range("E16").Value = '  -0.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '23.83'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.36%  '

$ws.Range("D18").Value = '2.490.13'
$ws.Range("E18").Value = '  +1.12%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.85'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.32%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '340.16'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.82'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.50%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.31%  '

$ws.Range("D27").Value = '2.604.42'
$ws.Range("E27").Value = '  +1.12%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.26%  '

$ws.Range("D30").Value = '0.0₃0829'
$ws.Range("E30").Value = '  -1.62%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '434.90'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.79%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.71'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.08'
$ws.Range("D37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.110'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.20%  '

$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.90'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.303'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.04%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.12%  '

$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '133.63'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.37'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0718'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.01%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.489'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("E50").Value = '  +0.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0918'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.39%  '

